$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.801.41'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').Value = '3.746.66'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '613.27'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '178.67'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.40%  '
$ws.Range('D7').Value = '3.748.48'
$ws.Range('E7').Value = '  +0.40%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.530'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.76%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.166'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.62%  '
$ws.Range('E11').Value = '  +3.69%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.483'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.95%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '40.13'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.23%  '
$ws.Range('E14').Value = '  +0.60%  '
$ws.Range('D15').Value = '4.366.06'
$ws.Range('E15').Value = '  +0.23%  '
$ws.Range('D16').Value = '3.743.49'
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('D17').Value = '69.795.56'
$ws.Range('E17').Value = '  +0.29%  '
$ws.Range('E18').Value = '  -2.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.46'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.39%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.42'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.59%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '502.69'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.18'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.47%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.722'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.62'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.70%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.09'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.10%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.37'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.95'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.41%  '
$ws.Range('E28').Value = '  +8.10%  '
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.48'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.00%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.13'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.58%  '
$ws.Range('E32').Value = '  +2.86%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '30.49'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.18%  '
$ws.Range('E34').Value = '  -1.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('E36').Value = '  +1.27%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.14'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.84%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.351'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.99%  '
$ws.Range('E39').Value = '  +3.76%  '
$ws.Range('B40').Value = 'dogwifhat'
$ws.Range('C40').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.09'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +14.05%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '446.16'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.56%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.07'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.86%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '49.73'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.80%  '
$ws.Range('B44').Value = 'Arweave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '45.31'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.57'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.41%  '
$ws.Range('D46').Value = '2.951.48'
$ws.Range('E46').Value = '  -4.45%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0360'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.77%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '27.16'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.26%  '
$ws.Range('B49').Value = 'USDe'
$ws.Range('C49').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.00'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '138.44'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.55%  '
$ws.Range('E51').Value = '  -0.75%  '
